# Apply the edit described by the diff:
# - Remove the "Real Madrid CF - Cadiz CF" / "19/12/2021" row (row 1),
#   shifting every following row up by one.
# - Change "Real Madrid CF - Valencia CF" date from 09/01/2022 to 08/01/2022.
# - Append a new row at the end: "Real Madrid CF - UD Levante" / "11/05/2022".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first data row (Real Madrid CF - Cadiz CF / 19/12/2021);
# this shifts all subsequent rows up by one, so the old row 2
# (Real Madrid CF - Valencia CF) becomes row 1, etc.
$ws.Range("A1:B1").Delete()

# Helper trick so date-like text ("dd/mm/yyyy") is stored as plain text
# instead of being auto-converted to a date serial number: write it as a
# formula that evaluates to the literal string, then collapse the formula
# down to its computed (text) value. This avoids touching any cell's
# number format / style.
$ws.Range("B1").Formula = '="08/01/2022"'
$ws.Range("B1").Value = $ws.Range("B1").Value

# Append the new fixture at the end of the table (row 10).
$ws.Range("A10").Value = "Real Madrid CF - UD Levante"
$ws.Range("B10").Formula = '="11/05/2022"'
$ws.Range("B10").Value = $ws.Range("B10").Value
